$errSubmit = "Error occurred during submission"
$elementFound = "Element Found"
$syntaxErr2 = "SyntaxError: bad input on line 2"
$creatingLinkedListTypo = "Creating Linked LIst"
$urlHeader = "url"
$url1 = "https://dsportalapp.herokuapp.com/question/1"
$url2 = "https://dsportalapp.herokuapp.com/question/2"
$url3 = "https://dsportalapp.herokuapp.com/question/3"
$url4 = "https://dsportalapp.herokuapp.com/question/4"
$submissionSuccess = "submission success"
$arrResult = "[4, 9, 9, 49, 121]"
$codeMaxConsec = "def findMaxConsecutiveOnes(nums) :`ncount = 0`nresult = 0`nfor i in range(0, len(nums)):`nif (nums[i] == 0):`ncount = 0`n`b`n`b`nelse:`ncount+= 1`n`b`n`b`nresult = max(result, count)`n`b`n`b`nprint(result)`n`b`n`b`nfindMaxConsecutiveOnes([1,0,1,1,0,1])"
$codeFindNumbers = "def findNumbers(nums):`nc=0`nfor i in nums:`nj=str(i)`nx=len(j)`nif x%2==0:`nc=c+1`n`b`n`b`n`b`n`b`nprint c`nreturn c`nfindNumbers([12,345,2,6,7896])"
$codeSortedSquares = "def sortedSquares(nums):`nsquares_list = []`nfor i in range(0, len(nums)):`nsquare = nums[i] * nums[i];`nsquares_list.append(square)`n`b`n`b`nsorted_squares_list = sorted(squares_list)`nprint sorted_squares_list;`nreturn sorted_squares_list;`nsortedSquares([-7,-3,2,3,11])"
$codeSearch = "def search(input_list, num):`nif(num in input_list):`nprint(`"Element Found`")`n`b`n`b`nelse:`nprint(`"Not Found`")`n`b`n`b`n`b`n`b`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$pythoncodeHeader = "pythoncode"
$two = "2"
$runHeader = "Run"
$submitHeader = "Submit"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet "login" (sheet1): a few new cells + a label correction
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item("login")

$login.Range("G3").Value = $errSubmit
$login.Range("G4").Value = $elementFound
$login.Range("B14").Value = $syntaxErr2
$login.Range("D16").Value = $creatingLinkedListTypo

# ---------------------------------------------------------------------------
# 2. Sheet "code" (sheet2): add a small list of question URLs + a hyperlink
# ---------------------------------------------------------------------------
$code = $wb.Worksheets.Item("code")

$code.Range("A1").Value = $urlHeader
$code.Range("A2").Value = $url1
$code.Range("A3").Value = $url2
$code.Range("A4").Value = $url3
$code.Range("A5").Value = $url4

$code.Range("A2").Font.Underline = $true
$code.Range("A2").Font.Color = 10498160
$code.Hyperlinks.Add($code.Range("A2"), $url1) | Out-Null

$code.Columns.Item(1).ColumnWidth = 44.5703125

# ---------------------------------------------------------------------------
# 3. New sheet "python": practice-question code samples + run/submit columns
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$python = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$python.Name = "python"

# widen columns first so autofit / wrap decisions below see the real widths
$python.Columns.Item(1).ColumnWidth = 243.42578125
$python.Columns.Item(2).ColumnWidth = 18.28515625
$python.Columns.Item(3).ColumnWidth = 18.28515625

# column B/C hold values that must stay plain TEXT ("2" as text, not 2 as a
# number) - set the text number format before writing the values
$python.Range("B2:C9").NumberFormat = "@"
$python.Range("A2:A9").WrapText = $false

$python.Range("A1").Value = $pythoncodeHeader
$python.Range("B1").Value = $runHeader
$python.Range("C1").Value = $submitHeader

$python.Range("A2").Value = $codeSearch
$python.Range("B2").Value = $elementFound
$python.Range("C2").Value = $submissionSuccess

$python.Range("A3").Value = $codeMaxConsec
$python.Range("B3").Value = $two
$python.Range("C3").Value = $submissionSuccess

$python.Range("A4").Value = $codeFindNumbers
$python.Range("B4").Value = $two
$python.Range("C4").Value = $submissionSuccess

$python.Range("A5").Value = $codeSortedSquares
$python.Range("B5").Value = $arrResult
$python.Range("C5").Value = $submissionSuccess

# header row formatting: white fill, black text, thin border, left aligned
$header = $python.Range("A1:C1")
$header.Interior.Pattern = 1
$header.Interior.Color = 16777215
$header.Interior.PatternColor = 16777215
$header.Font.Color = 0
$header.Font.Name = "Calibri"
$header.Font.Size = 11
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.HorizontalAlignment = -4131

# body rows formatting: thin border around every used cell
$body = $python.Range("A2:C9")
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2
$body.Font.Color = 0
$body.Font.Name = "Calibri"
$body.Font.Size = 11

# the multi-line code in column A pushed auto row-height up; since the
# column is wide enough that the text never actually wraps, restore every
# data row back to the sheet's natural (default) height
for ($r = 2; $r -le 9; $r++) {
    $python.Rows.Item($r).AutoFit()
}

$python.Range("A4").Select() | Out-Null
